$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 5.010563347124673
$ws.Range("C2").Value = 0.6554134751303877
$ws.Range("D2").Value = 0.07579977628580536
$ws.Range("E2").Value = 0.02446419001990918
$ws.Range("G2").Value = 0.002655987569348779
$ws.Range("L2").Value = 0.1988097257590979
$ws.Range("M2").Value = 0.7849432847482802
$ws.Range("N2").Value = 4.316152919285202
$ws.Range("B3").Value = 4.825516395299758
$ws.Range("C3").Value = 0.5992126970652407
$ws.Range("D3").Value = 0.06915090357891529
$ws.Range("E3").Value = 0.02205428770693274
$ws.Range("G3").Value = 0.002665069102175239
$ws.Range("L3").Value = 0.1968057160936851
$ws.Range("M3").Value = 0.7605971800985074
$ws.Range("N3").Value = 4.233266226938355
$ws.Range("B4").Value = 4.715444826874659
$ws.Range("C4").Value = 0.5651654514105076
$ws.Range("D4").Value = 0.06512091509036111
$ws.Range("E4").Value = 0.02056768052905511
$ws.Range("G4").Value = 0.002670924753622492
$ws.Range("L4").Value = 0.195685027070823
$ws.Range("M4").Value = 0.746188586931865
$ws.Range("N4").Value = 4.182794913455552
$ws.Range("B5").Value = 4.67147095682077
$ws.Range("C5").Value = 0.5514033218601071
$ws.Range("D5").Value = 0.06349146509153059
$ws.Range("E5").Value = 0.01995998826416567
$ws.Range("G5").Value = 0.002673381585602196
$ws.Range("L5").Value = 0.1952558168489134
$ws.Range("M5").Value = 0.7404513102987096
$ws.Range("N5").Value = 4.162329826940095
$ws.Range("B6").Value = 4.664222034359625
$ws.Range("C6").Value = 0.5491248230547967
$ws.Range("D6").Value = 0.06322165831453219
$ws.Range("E6").Value = 0.0198589631177768
$ws.Range("G6").Value = 0.002673793814062061
$ws.Range("L6").Value = 0.1951862029313958
$ws.Range("M6").Value = 0.7395067121248715
$ws.Range("N6").Value = 4.158937694407769
$ws.Range("B7").Value = 4.714848226312654
$ws.Range("C7").Value = 0.5649793997803272
$ws.Range("D7").Value = 0.06509888842252565
$ws.Range("E7").Value = 0.02055949281677272
$ws.Range("G7").Value = 0.00267095760115994
$ws.Range("L7").Value = 0.1956791274898677
$ws.Range("M7").Value = 0.7461106699538078
$ws.Range("N7").Value = 4.182518504314487
$ws.Range("B8").Value = 4.946016016304782
$ws.Range("C8").Value = 0.6359378764357757
$ws.Range("D8").Value = 0.07349612394546057
$ws.Range("E8").Value = 0.02363460086198899
$ws.Range("G8").Value = 0.002659061048153094
$ws.Range("L8").Value = 0.1980958801930868
$ws.Range("M8").Value = 0.7764357769334893
$ws.Range("N8").Value = 4.28748391631072
$ws.Range("B9").Value = 5.428021407016445
$ws.Range("C9").Value = 0.7788992176299985
$ws.Range("D9").Value = 0.09039769157806177
$ws.Range("E9").Value = 0.02961700387477251
$ws.Range("G9").Value = 0.002637935849428428
$ws.Range("L9").Value = 0.2037125302050171
$ws.Range("M9").Value = 0.8402586661214571
$ws.Range("N9").Value = 4.496838334729034
$ws.Range("B10").Value = 5.800434950874092
$ws.Range("C10").Value = 0.8864924200719884
$ws.Range("D10").Value = 0.1031074425154657
$ws.Range("E10").Value = 0.03399327314641454
$ws.Range("G10").Value = 0.002623738701323118
$ws.Range("L10").Value = 0.2083837435958742
$ws.Range("M10").Value = 0.8899102479541341
$ws.Range("N10").Value = 4.65305245163313
$ws.Range("B11").Value = 5.974006078039338
$ws.Range("C11").Value = 0.9360481020939915
$ws.Range("D11").Value = 0.1089589956685444
$ws.Range("E11").Value = 0.03598229951217391
$ws.Range("G11").Value = 0.002617563118453365
$ws.Range("L11").Value = 0.2106293218228359
$ws.Range("M11").Value = 0.9131215206894154
$ws.Range("N11").Value = 4.72469708499267
$ws.Range("B12").Value = 6.040345261771336
$ws.Range("C12").Value = 0.9549056414539905
$ws.Range("D12").Value = 0.111185354805329
$ws.Range("E12").Value = 0.03673542185122969
$ws.Range("G12").Value = 0.002615264910305925
$ws.Range("L12").Value = 0.2114971869790168
$ws.Range("M12").Value = 0.9220027372547293
$ws.Range("N12").Value = 4.751915187933491
$ws.Range("B13").Value = 6.026030511242993
$ws.Range("C13").Value = 0.9508401815588172
$ws.Range("D13").Value = 0.1107053939208384
$ws.Range("E13").Value = 0.03657322416414033
$ws.Range("G13").Value = 0.002615758080883729
$ws.Range("L13").Value = 0.2113094950270238
$ws.Range("M13").Value = 0.9200859045810006
$ws.Range("N13").Value = 4.746049309196223
$ws.Range("B14").Value = 5.979451506617124
$ws.Range("C14").Value = 0.9375976584377099
$ws.Range("D14").Value = 0.109141946387652
$ws.Range("E14").Value = 0.03604425990513249
$ws.Range("G14").Value = 0.002617373236683057
$ws.Range("L14").Value = 0.2107003696290093
$ws.Range("M14").Value = 0.9138503371544573
$ws.Range("N14").Value = 4.726934550376313
$ws.Range("B15").Value = 5.951000575589546
$ws.Range("C15").Value = 0.9294983169527882
$ws.Range("D15").Value = 0.1081856713542209
$ws.Range("E15").Value = 0.03572024901077597
$ws.Range("G15").Value = 0.002618367811905311
$ws.Range("L15").Value = 0.2103295486064383
$ws.Range("M15").Value = 0.9100428570898202
$ws.Range("N15").Value = 4.715237771553348
$ws.Range("B16").Value = 5.789176499536893
$ws.Range("C16").Value = 0.883266475546975
$ws.Range("D16").Value = 0.1027264740436493
$ws.Range("E16").Value = 0.03386326498274173
$ws.Range("G16").Value = 0.002624147953679972
$ws.Range("L16").Value = 0.2082394312122631
$ws.Range("M16").Value = 0.8884060567887957
$ws.Range("N16").Value = 4.648382370490538
$ws.Range("B17").Value = 5.690977273323483
$ws.Range("C17").Value = 0.8550640544568182
$ws.Range("D17").Value = 0.09939564211136087
$ws.Range("E17").Value = 0.03272373730688471
$ws.Range("G17").Value = 0.002627766086752083
$ws.Range("L17").Value = 0.2069882233843146
$ws.Range("M17").Value = 0.8752938031570778
$ws.Range("N17").Value = 4.607520742480659
$ws.Range("B18").Value = 5.634886240977949
$ws.Range("C18").Value = 0.8389000308749246
$ws.Range("D18").Value = 0.09748637944960592
$ws.Range("E18").Value = 0.03206813474575654
$ws.Range("G18").Value = 0.002629873775707739
$ws.Range("L18").Value = 0.2062799045912556
$ws.Range("M18").Value = 0.8678106531074263
$ws.Range("N18").Value = 4.584072747765731
$ws.Range("B19").Value = 5.615961453148316
$ws.Range("C19").Value = 0.8334368845258382
$ws.Range("D19").Value = 0.0968410459250606
$ws.Range("E19").Value = 0.03184612344508153
$ws.Range("G19").Value = 0.002630591986883285
$ws.Range("L19").Value = 0.2060420231355096
$ws.Range("M19").Value = 0.8652870097901513
$ws.Range("N19").Value = 4.576142900924594
$ws.Range("B20").Value = 5.70139022931312
$ws.Range("C20").Value = 0.8580602950408434
$ws.Range("D20").Value = 0.09974953428982758
$ws.Range("E20").Value = 0.03284505897787326
$ws.Range("G20").Value = 0.002627378175903084
$ws.Range("L20").Value = 0.2071202414509656
$ws.Range("M20").Value = 0.8766835397587727
$ws.Range("N20").Value = 4.611864858872536
$ws.Range("B21").Value = 5.993116181298433
$ws.Range("C21").Value = 0.941484782395662
$ws.Range("D21").Value = 0.1096008800360693
$ws.Range("E21").Value = 0.03619963025818507
$ws.Range("G21").Value = 0.002616897733810156
$ws.Range("L21").Value = 0.2108788078412971
$ws.Range("M21").Value = 0.9156793734900361
$ws.Range("N21").Value = 4.732546601187096
$ws.Range("B22").Value = 6.187347221678579
$ws.Range("C22").Value = 0.9965449840530027
$ws.Range("D22").Value = 0.1161007464401962
$ws.Range("E22").Value = 0.03839166126097027
$ws.Range("G22").Value = 0.002610283208415165
$ws.Range("L22").Value = 0.2134374125730574
$ws.Range("M22").Value = 0.9417002633620939
$ws.Range("N22").Value = 4.811933045080195
$ws.Range("B23").Value = 6.083351090663882
$ws.Range("C23").Value = 0.9671077838645488
$ws.Range("D23").Value = 0.1126258675111842
$ws.Range("E23").Value = 0.03722170869183117
$ws.Range("G23").Value = 0.002613792101288259
$ws.Range("L23").Value = 0.2120624317238224
$ws.Range("M23").Value = 0.927762874141564
$ws.Range("N23").Value = 4.769514591306745
$ws.Range("B24").Value = 5.696681396006738
$ws.Range("C24").Value = 0.8567055394736371
$ws.Range("D24").Value = 0.09958952198847726
$ws.Range("E24").Value = 0.03279021093663204
$ws.Range("G24").Value = 0.002627553464365907
$ws.Range("L24").Value = 0.2070605217796384
$ws.Range("M24").Value = 0.8760550678568251
$ws.Range("N24").Value = 4.609900747426593
$ws.Range("B25").Value = 5.294474093087501
$ws.Range("C25").Value = 0.7397922683657043
$ws.Range("D25").Value = 0.08577606596887222
$ws.Range("E25").Value = 0.02800269108135822
$ws.Range("G25").Value = 0.002643416898370291
$ws.Range("L25").Value = 0.2020981558573069
$ws.Range("M25").Value = 0.8225158504904186
$ws.Range("N25").Value = 4.439799518336486
